$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 3) down to the new row 4
# so the new row inherits the same cell styles (e.g. date format on column A).
$ws.Range("A3:N3").Copy() | Out-Null
$ws.Range("A4:N4").PasteSpecial(-4122) | Out-Null

$row = 4
$ws.Cells.Item($row, 1).Value = 42605.886712962965
$ws.Cells.Item($row, 2).Value = -24
$ws.Cells.Item($row, 3).Value = 50
$ws.Cells.Item($row, 4).Value = 48
$ws.Cells.Item($row, 5).Value = 22
$ws.Cells.Item($row, 6).Value = 77
$ws.Cells.Item($row, 7).Value = 17007
$ws.Cells.Item($row, 8).Value = 5713
$ws.Cells.Item($row, 9).Value = 1088
$ws.Cells.Item($row, 10).Value = 105
$ws.Cells.Item($row, 11).Value = 101
$ws.Cells.Item($row, 12).Value = 6
$ws.Cells.Item($row, 13).Value = 21
$ws.Cells.Item($row, 14).Value = "Noun"
